# Penalty Reward System tweak (unfinished) -- shifts the forecast week
# dates forward by one week and updates MyForecast figures, then
# refreshes the derived metrics on the Summary sheet to match.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# New Week_Start_Date (col B) and MyForecast (col D) values per row.
$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$newForecast = @(125, 137, 89, 124, 113, 107, 106, 102, 95, 86, 83, 85, 89, 91, 87, 89)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $bCell = $wsForecast.Range("B$row")
    # Keep these as literal text (not auto-converted to date serials).
    $bCell.NumberFormat = "@"
    $bCell.Value = $newDates[$i]

    $wsForecast.Range("D$row").Value = $newForecast[$i]
}

# --- Summary sheet --------------------------------------------------------
# "2023-01-01 to 2025-01-05" is not date/number-like on its own, so it is
# safely kept as text without forcing a number format.
$wsSummary.Range("B2").Value = "2023-01-01 to 2025-01-05"

$b4 = $wsSummary.Range("B4")
$b4.NumberFormat = "@"
$b4.Value = "208"

$b6 = $wsSummary.Range("B6")
$b6.NumberFormat = "@"
$b6.Value = "79"

$wsSummary.Range("B8").Value = "8806 units"

$b9 = $wsSummary.Range("B9")
$b9.NumberFormat = "@"
$b9.Value = "1608"

$b10 = $wsSummary.Range("B10")
$b10.NumberFormat = "@"
$b10.Value = "903"

$b11 = $wsSummary.Range("B11")
$b11.NumberFormat = "@"
$b11.Value = "475"

$b12 = $wsSummary.Range("B12")
$b12.NumberFormat = "@"
$b12.Value = "137"

$b14 = $wsSummary.Range("B14")
$b14.NumberFormat = "@"
$b14.Value = "83"

$b15 = $wsSummary.Range("B15")
$b15.NumberFormat = "@"
$b15.Value = "2025-03-23"
